# Apply the TUESDAY schedule updates to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AV2").Value = "RLN"

$ws.Range("L5").Value = "DK"
$ws.Range("P5").Value = "TT"
$ws.Range("T5").Value = "TT"
$ws.Range("AJ5").Value = "DK"
$ws.Range("AN5").Value = "TT"

$ws.Range("AB6").Value = "AMS"
$ws.Range("AF6").Value = "AMS"

$ws.Range("AJ7").Value = "ZL"

$ws.Range("AB8").Value = "AS/SH"
$ws.Range("AF8").Value = "FD/BK"

$ws.Range("P11").Value = "MP"
$ws.Range("T11").Value = "MP"
$ws.Range("AN11").Value = "MP"

$ws.Range("AC24").Value = 62

$ws.Range("AB25").Value = "FGN/CK"
$ws.Range("AF25").Value = "FGN/CK"
